$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelle photo pour Francois Mercier : jpeg -> png
# (leading apostrophe preserves the text/quote-prefix cell style, s="2",
#  exactly as the avatar column cells already use)
$ws.Range("C5").Formula = "'/assets/francois-mercier.png"

# Move the active selection to C6
$ws.Range("C6").Select()
